# Replace the MS-Word field-code representation of the M2Doc tags
# ( {m:for ...}, {m:v.name}, {m:endfor} ) with plain literal text runs
# containing the same braces/tags, as produced by
# TokenIteratorFieldRewriterSplit.
#
# Each target paragraph keeps its own <w:p> (and any runs that are not
# part of the field) untouched; only the fldChar begin/instrText.../fldChar
# end run sequence is rewritten into plain <w:t> runs.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 2 : "{m:for v | self.eClassifiers}" ------------------------
$p2 = $d.Paragraphs.Item(2)
$xml2 = "<w:p $wNs w:rsidP='00F5495F' w:rsidR='00052FB8' w:rsidRDefault='00C52979'>" +
        "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
        "<w:r><w:t>{m:</w:t></w:r>" +
        "<w:r><w:t>f</w:t></w:r>" +
        "<w:r><w:t>r v | self.eClassifiers}</w:t></w:r>" +
        "</w:p>"
[void]$p2.Range.InsertXML($xml2)

# --- Paragraph 3 : "name<nbsp>= {m:v.name}," -------------------------------
$p3 = $d.Paragraphs.Item(3)
$nbsp = [char]0x00A0
$xml3 = "<w:p $wNs w:rsidP='00F5495F' w:rsidR='00052FB8' w:rsidRDefault='00730F02'>" +
        "<w:proofErr w:type='spellStart'/>" +
        "<w:r><w:t>name</w:t></w:r>" +
        "<w:proofErr w:type='spellEnd'/>" +
        "<w:r><w:t>$nbsp</w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:t>=</w:t></w:r>" +
        "<w:r w:rsidR='00E27251'><w:t xml:space='preserve'> </w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:t>{</w:t></w:r>" +
        "<w:r><w:t>m</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>:v.name}</w:t></w:r>" +
        "<w:r w:rsidR='00052FB8'><w:t>,</w:t></w:r>" +
        "</w:p>"
[void]$p3.Range.InsertXML($xml3)

# --- Paragraph 4 : "{m:endfor}" --------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$xml4 = "<w:p $wNs w:rsidP='00F5495F' w:rsidR='00C52979' w:rsidRDefault='006F5523'>" +
        "<w:r><w:t>{</w:t></w:r>" +
        "<w:r><w:t>m:</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>endfor}</w:t></w:r>" +
        "</w:p>"
[void]$p4.Range.InsertXML($xml4)
